$wb = $excel.ActiveWorkbook

# Add a new worksheet after the last existing sheet (becomes active) and
# name it "Sheet2" -- mirrors the author re-running the simulation and
# pasting a fresh trace of the orderbook/matching output onto a new tab.
$count = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Sheet2"

# --- Orderbook / trade trace table -----------------------------------
# Entered in the same order the author typed them (first the bid that
# got matched twice by mistake, then the asks, then the remaining bids,
# and finally the last partial-fill ask row) so the shared-string table
# comes out in the same sequence as the authored workbook.

# Row 8 (bid, first entered)
$ws.Range("B8").Value = "T(uuid: 0"
$ws.Range("C8").Value = " order id: 1"
$ws.Range("D8").Value = " Type: OrderType.BID"
$ws.Range("E8").Value = " Volume: 20"
$ws.Range("F8").Value = " Price: 49.5)"
$ws.Range("G8").Value = " "

# Row 3 (ask)
$ws.Range("B3").Value = "T(uuid: 1"
$ws.Range("C3").Value = " order id: 0"
$ws.Range("D3").Value = " Type: OrderType.ASK"
$ws.Range("E3").Value = " Volume: 20"
$ws.Range("F3").Value = " Price: 49.5)"
$ws.Range("G3").Value = " "

# Row 4 (ask)
$ws.Range("B4").Value = "T(uuid: 2"
$ws.Range("C4").Value = " order id: 0"
$ws.Range("D4").Value = " Type: OrderType.ASK"
$ws.Range("E4").Value = " Volume: 50"
$ws.Range("F4").Value = " Price: 39.5)"
$ws.Range("G4").Value = " "

# Row 5 (ask)
$ws.Range("B5").Value = "T(uuid: 4"
$ws.Range("C5").Value = " order id: 0"
$ws.Range("D5").Value = " Type: OrderType.ASK"
$ws.Range("E5").Value = " Volume: 42"
$ws.Range("F5").Value = " Price: 33.0)"
$ws.Range("G5").Value = " "

# Row 6 (ask)
$ws.Range("B6").Value = "T(uuid: 3"
$ws.Range("C6").Value = " order id: 0"
$ws.Range("D6").Value = " Type: OrderType.ASK"
$ws.Range("E6").Value = " Volume: 38"
$ws.Range("F6").Value = " Price: 33.0)"
$ws.Range("G6").Value = " "

# Row 9 (bid)
$ws.Range("B9").Value = "T(uuid: 0"
$ws.Range("C9").Value = " order id: 1"
$ws.Range("D9").Value = " Type: OrderType.BID"
$ws.Range("E9").Value = " Volume: 50"
$ws.Range("F9").Value = " Price: 39.5)"
$ws.Range("G9").Value = " "

# Row 10 (bid)
$ws.Range("B10").Value = "T(uuid: 0"
$ws.Range("C10").Value = " order id: 4"
$ws.Range("D10").Value = " Type: OrderType.BID"
$ws.Range("E10").Value = " Volume: 33"
$ws.Range("F10").Value = " Price: 33.0)"
$ws.Range("G10").Value = " "

# Row 11 (bid)
$ws.Range("B11").Value = "T(uuid: 0"
$ws.Range("C11").Value = " order id: 3"
$ws.Range("D11").Value = " Type: OrderType.BID"
$ws.Range("E11").Value = " Volume: 5"
$ws.Range("F11").Value = " Price: 33.5)"
$ws.Range("G11").Value = " "

# Row 7 (ask, last partial fill, added last)
$ws.Range("B7").Value = "T(uuid: 4"
$ws.Range("C7").Value = " order id: 0"
$ws.Range("D7").Value = " Type: OrderType.ASK"
$ws.Range("E7").Value = " Volume: 3"
$ws.Range("F7").Value = " Price: 33.5)]"

# --- Volume sanity-check block (rows 16-20) ---------------------------
$ws.Range("C16").Value = " Volume"
$ws.Range("D16").Value = 20
$ws.Range("E16").Formula = "=D16"

$ws.Range("C17").Value = " Volume"
$ws.Range("D17").Value = 50
$ws.Range("E17").Formula = "=D17"

$ws.Range("C18").Value = " Volume"
$ws.Range("D18").Value = 42
$ws.Range("E18").Formula = "=D18"

$ws.Range("C19").Value = " Volume"
$ws.Range("D19").Value = 38
$ws.Range("E19").Formula = "=D19"

$ws.Range("C20").Value = " Volume"
$ws.Range("D20").Value = 3
$ws.Range("E20").Formula = "=D20"

# --- Net exposure block (rows 24-27), volumes negated -----------------
$ws.Range("C24").Value = " Volume"
$ws.Range("D24").Value = 20
$ws.Range("E24").Formula = "=0-D24"

$ws.Range("C25").Value = " Volume"
$ws.Range("D25").Value = 50
$ws.Range("E25").Formula = "=0-D25"

$ws.Range("C26").Value = " Volume"
$ws.Range("D26").Value = 33
$ws.Range("E26").Formula = "=0-D26"

$ws.Range("C27").Value = " Volume"
$ws.Range("D27").Value = 5
$ws.Range("E27").Formula = "=0-D27"

# --- Column widths (as authored) --------------------------------------
$ws.Columns("B").ColumnWidth = 39.7109375
$ws.Columns("C").ColumnWidth = 35.5703125
$ws.Columns("D").ColumnWidth = 40.5703125
$ws.Columns("E").ColumnWidth = 50
$ws.Columns("F").ColumnWidth = 39.140625

# --- Sort the trade table by volume (column D), matching the recorded sortState ---
$ws.Range("B3:G11").Sort($ws.Range("D3:D11"), 1)

# --- Final selection on the new sheet ---
$ws.Range("F14").Select()
